# Applies crypto price/volume updates from the Jan 1 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch keep their original "Text" storage type so that
# values such as "1.00" or "3.77" are not silently reinterpreted as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.802.89"
$ws.Range("E2").Value = "  +2.13%  "

$ws.Range("D3").Value = "2.335.29"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "311.48"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").Value = "108.34"
$ws.Range("E6").Value = "  +2.99%  "

$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  +1.92%  "

$ws.Range("D10").Value = "41.17"
$ws.Range("E10").Value = "  +3.75%  "

$ws.Range("D11").Value = "0.0917"
$ws.Range("E11").Value = "  +1.39%  "

$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").Value = "1.00"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").Value = "15.44"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").Value = "2.692.78"
$ws.Range("E16").Value = "  +2.18%  "

$ws.Range("D17").Value = "2.333.09"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").Value = "43.755.18"
$ws.Range("E18").Value = "  +2.23%  "

$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").Value = "13.05"
$ws.Range("E21").Value = "  -5.25%  "

$ws.Range("D22").Value = "74.18"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("E23").Value = "  -3.62%  "

$ws.Range("D24").Value = "268.75"
$ws.Range("E24").Value = "  +1.20%  "

$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").Value = "7.68"
$ws.Range("E27").Value = "  +6.65%  "

$ws.Range("D28").Value = "11.11"
$ws.Range("E28").Value = "  +2.54%  "

$ws.Range("E29").Value = "  -1.85%  "

$ws.Range("D30").Value = "39.05"
$ws.Range("E30").Value = "  +4.76%  "

$ws.Range("D31").Value = "22.57"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").Value = "168.62"
$ws.Range("E32").Value = "  +0.78%  "

$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("E34").Value = "  +9.37%  "

$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("D37").Value = "4.72"
$ws.Range("E37").Value = "  +3.33%  "

$ws.Range("D38").Value = "0.0363"
$ws.Range("E38").Value = "  +3.02%  "

$ws.Range("E39").Value = "  +8.10%  "

$ws.Range("D40").Value = "3.77"

$ws.Range("E41").Value = "  +7.99%  "

$ws.Range("D42").Value = "105.13"
$ws.Range("E42").Value = "  +11.68%  "

$ws.Range("E43").Value = "  +2.28%  "

$ws.Range("D44").Value = "13.39"
$ws.Range("E44").Value = "  +10.12%  "

$ws.Range("D45").Value = "71.54"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").Value = "114.13"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").Value = "1.672.76"
$ws.Range("E48").Value = "  -3.83%  "

$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "0.218"
$ws.Range("E49").Value = "  +14.72%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "76.73"
$ws.Range("E50").Value = "  -4.22%  "

$ws.Range("E51").Value = "  +2.01%  "
